# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3140
$ws1.Range("F5").Value = 6847
$ws1.Range("F6").Value = 1906
$ws1.Range("F7").Value = 11
$ws1.Range("F8").Value = 67
$ws1.Range("F12").Value = 17
$ws1.Range("F13").Value = 144
$ws1.Range("F14").Value = 164

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3140
$ws4.Range("F6").Value = 6847
$ws4.Range("F7").Value = 1906
$ws4.Range("F8").Value = 11
$ws4.Range("F9").Value = 67
$ws4.Range("F13").Value = 17
$ws4.Range("F14").Value = 144
$ws4.Range("F15").Value = 164
